# Apply changes described by the commit: "Added new patient for LAPhases, corrected openfcn"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Correct a rounding artefact on AI3 (SystoleLength/DiastoleLength)
# ---------------------------------------------------------------------
$ws.Range("AI3").Value = 0.419279907084785

# ---------------------------------------------------------------------
# 2) Row 56 used to hold a dummy "AristotelesTest" row - repurpose it to
#    describe the LBBB test patient, keeping only the MVO/MVC/AVO/AVC
#    inputs and clearing everything that depended on the ECG curve
#    selection (it is not available for this patient).
# ---------------------------------------------------------------------
$ws.Range("A56").Value = "LBBB_Teste"
$ws.Range("B56").Value = "LBBB"
$ws.Range("Q56").Value = 581
$ws.Range("R56").Value = 33
$ws.Range("S56").Value = 176
$ws.Range("T56").Value = 467
$ws.Range("U56:AK56").ClearContents()

# ---------------------------------------------------------------------
# 3) Add a new patient (NataliaOno) in row 57, reusing the formatting
#    that row 56 already has (so the new row looks consistent with the
#    rest of the table).
# ---------------------------------------------------------------------
$ws.Range("A56:AK56").Copy()
$ws.Range("A57:AK57").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A57").Value = "NataliaOno"
$ws.Range("Q57").Value = 437
$ws.Range("R57").Value = 13
$ws.Range("S57").Value = 36
$ws.Range("T57").Value = 367
$ws.Range("U57").Value = 169
$ws.Range("V57").Value = 1134
$ws.Range("W57").Value = 1325
$ws.Range("X57").Value = 169
$ws.Range("Y57").Value = 240
$ws.Range("Z57").Value = 270
$ws.Range("AA57").Value = 597
$ws.Range("AB57").Value = 670
$ws.Range("AC57").Value = 1134
$ws.Range("AD57").Value = 1325
$ws.Range("AE57").Value = 1620
$ws.Range("AF57").Value = 1650
$ws.Range("AG57").Value = 357
$ws.Range("AH57").Value = 1023
$ws.Range("AI57").Value = 0.3489736070381232
$ws.Range("AJ57").Value = -20.9
$ws.Range("AK57").Value = 30.7

# ---------------------------------------------------------------------
# 4) Update the view so the new row is the active cell / visible area,
#    matching what the author had selected when saving the workbook.
# ---------------------------------------------------------------------
$aw = $excel.ActiveWindow
$aw.ScrollRow = 16
$aw.ScrollColumn = 1
$ws.Range("A57").Select()
